$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format ("@") on the Price cells we are about to overwrite,
# one at a time, so Excel stores the literal scraped string (e.g. "310.02",
# "19.40", "2.517.34") instead of silently re-parsing it as a float and losing
# trailing zeros / the thousands-dot formatting used throughout this column.
$priceCells = @("D2", "D3", "D5", "D6", "D9", "D10", "D13", "D14", "D15", "D16", "D18", "D19", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D30", "D31", "D33", "D35", "D38", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($c in $priceCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value2 = "43.206.55"
$ws.Range("E2").Value2 = "  +2.11%  "

$ws.Range("D3").Value2 = "2.295.16"
$ws.Range("E3").Value2 = "  +0.89%  "

$ws.Range("E4").Value2 = "  +0.00%  "

$ws.Range("D5").Value2 = "310.02"
$ws.Range("E5").Value2 = "  +1.29%  "

$ws.Range("D6").Value2 = "101.19"
$ws.Range("E6").Value2 = "  +3.88%  "

$ws.Range("E7").Value2 = "  +0.18%  "

$ws.Range("E8").Value2 = "  +0.05%  "

$ws.Range("D9").Value2 = "0.507"
$ws.Range("E9").Value2 = "  +2.68%  "

$ws.Range("D10").Value2 = "36.28"
$ws.Range("E10").Value2 = "  +1.10%  "

$ws.Range("E11").Value2 = "  +2.62%  "

$ws.Range("E12").Value2 = "  +0.42%  "

$ws.Range("D13").Value2 = "6.94"
$ws.Range("E13").Value2 = "  +3.95%  "

$ws.Range("D14").Value2 = "2.648.51"
$ws.Range("E14").Value2 = "  +1.76%  "

$ws.Range("D15").Value2 = "14.89"
$ws.Range("E15").Value2 = "  +2.86%  "

$ws.Range("D16").Value2 = "2.298.43"
$ws.Range("E16").Value2 = "  +1.55%  "

$ws.Range("D18").Value2 = "43.142.33"
$ws.Range("E18").Value2 = "  +2.19%  "

$ws.Range("D19").Value2 = "12.59"
$ws.Range("E19").Value2 = "  +0.50%  "

$ws.Range("E20").Value2 = "  +0.63%  "

$ws.Range("E21").Value2 = "  +0.80%  "

$ws.Range("D22").Value2 = "67.94"
$ws.Range("E22").Value2 = "  +0.04%  "

$ws.Range("D23").Value2 = "240.72"
$ws.Range("E23").Value2 = "  +0.70%  "

$ws.Range("D24").Value2 = "2.02"
$ws.Range("E24").Value2 = "  +3.15%  "

$ws.Range("D25").Value2 = "2.62"
$ws.Range("E25").Value2 = "  +1.01%  "

$ws.Range("D26").Value2 = "1.01"
$ws.Range("E26").Value2 = "  +1.25%  "

$ws.Range("D27").Value2 = "38.90"
$ws.Range("E27").Value2 = "  +3.89%  "

$ws.Range("D28").Value2 = "23.91"
$ws.Range("E28").Value2 = "  +0.44%  "

$ws.Range("E29").Value2 = "  +0.97%  "

$ws.Range("D30").Value2 = "2.16"
$ws.Range("E30").Value2 = "  +2.13%  "

$ws.Range("D31").Value2 = "165.78"
$ws.Range("E31").Value2 = "  +3.95%  "

$ws.Range("E32").Value2 = "  +0.61%  "

$ws.Range("D33").Value2 = "1.00"
$ws.Range("E33").Value2 = "  +0.04%  "

$ws.Range("E34").Value2 = "  -1.85%  "

$ws.Range("D35").Value2 = "17.98"
$ws.Range("E35").Value2 = "  +3.42%  "

$ws.Range("E36").Value2 = "  -0.39%  "

$ws.Range("E37").Value2 = "  +1.04%  "

$ws.Range("D38").Value2 = "0.106"
$ws.Range("E38").Value2 = "  -0.23%  "

$ws.Range("E39").Value2 = "  -0.59%  "

$ws.Range("E40").Value2 = "  +0.36%  "

$ws.Range("E41").Value2 = "  +2.00%  "

$ws.Range("E42").Value2 = "  -5.89%  "

$ws.Range("D43").Value2 = "19.40"
$ws.Range("E43").Value2 = "  +1.94%  "

$ws.Range("D44").Value2 = "0.0290"
$ws.Range("E44").Value2 = "  +1.52%  "

$ws.Range("D45").Value2 = "1.957.09"
$ws.Range("E45").Value2 = "  -1.68%  "

$ws.Range("D46").Value2 = "3.02"
$ws.Range("E46").Value2 = "  +2.57%  "

$ws.Range("D47").Value2 = "9.86"
$ws.Range("E47").Value2 = "  -0.89%  "

$ws.Range("D48").Value2 = "55.12"
$ws.Range("E48").Value2 = "  +3.27%  "

$ws.Range("E51").Value2 = "  +1.10%  "

# Rows 49 and 50 swapped rank order: RocketPoolETH moved above Stacks,
# both rows also carry refreshed Price / Volume(1h) figures.
$ws.Range("B49").Value2 = "RocketPoolETH"
$ws.Range("C49").Value2 = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value2 = "2.517.34"
$ws.Range("E49").Value2 = "  +1.06%  "

$ws.Range("B50").Value2 = "Stacks"
$ws.Range("C50").Value2 = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value2 = "1.53"
$ws.Range("E50").Value2 = "  +0.65%  "
